$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item("Сводный").Name = "Consolidated budget"
$wb.Worksheets.Item("Продажи").Name = "Sales"
$wb.Worksheets.Item("Производство").Name = "Production"
$wb.Worksheets.Item("Логистика").Name = "Logistics"
$wb.Worksheets.Item("Прочее").Name = "Misc"

$wsCB = $wb.Worksheets.Item("Consolidated budget")
$wsCB.Activate()
$wsCB.Range("E10").Select()

$wsSales = $wb.Worksheets.Item("Sales")
$wsSales.Activate()
$wsSales.Range("C34").Select()

$wsProd = $wb.Worksheets.Item("Production")
$wsProd.Activate()
$wsProd.Range("E34").Select()

$wsLog = $wb.Worksheets.Item("Logistics")
$wsLog.Activate()
$wsLog.Range("F33").Select()

$wsMisc = $wb.Worksheets.Item("Misc")
$wsMisc.Activate()
$wsMisc.Range("H34").Select()

$wsCB.Activate()
